# Daily attendance processing - 2026-01-04 09:56:54
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" for every row where it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
